$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (the second paragraph in the
#    document, directly after the title) in its entirety.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Annihilator Slot Game Free - Review &
#    Features" right before the final paragraph (the one that used to hold
#    the italic image-prompt text).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newRange = $newPara.Range
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Annihilator Slot Game Free - Review &amp; Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($xml)

# 3. Replace the text of the final (italic) paragraph with the review blurb,
#    keeping its existing italic run formatting intact.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$textOnly = $d.Range($finalRange.Start, $finalRange.End - 1)
$textOnly.Text = "Read our review of the Annihilator slot game. Play for free and find out about the game's features, including expanding wilds and an immersive soundtrack."
